$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.647.90"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.40%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.437.86"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.71%  "

# Row 4
$ws.Range("E4").Value = "  +0.10%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.57%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.19%  "

# Row 7
$ws.Range("E7").Value = "  -0.10%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.533"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.88%  "

# Row 9
$ws.Range("E9").Value = "  +2.82%  "

# Row 10
$ws.Range("E10").Value = "  +0.53%  "

# Row 11
$ws.Range("E11").Value = "  +2.27%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.356"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.79%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "26.88"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.84%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000181"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.60%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.889.93"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.32%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.565.55"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.43%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.438.02"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.94%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.27"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.57%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.98"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.27%  "

# Row 20
$ws.Range("E20").Value = "  +1.58%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.18"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.25%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.05%  "

# Row 23
$ws.Range("E23").Value = "  +7.75%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.33"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.33%  "

# Row 25
$ws.Range("B25").Value = "Bittensor"
$ws.Range("C25").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "585.73"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.52%  "

# Row 26
$ws.Range("B26").Value = "Aptos"
$ws.Range("C26").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.58"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.92%  "

# Row 27
$ws.Range("E27").Value = "  +9.79%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.556.43"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.64%  "

# Row 29
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.45"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.50%  "

# Row 30
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.73%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.45"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.83%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.146"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.41%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.88"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.91%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.50"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.14%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.85"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.29%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.998"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.15%  "

# Row 37
$ws.Range("E37").Value = "  +1.39%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.77"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.05%  "

# Row 39
$ws.Range("E39").Value = "  +0.30%  "

# Row 40
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.84"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.76%  "

# Row 41
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "148.07"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.64%  "

# Row 42
$ws.Range("E42").Value = "  +0.13%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.45"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +11.35%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "148.57"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.51%  "

# Row 45
$ws.Range("E45").Value = "  +3.07%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0535"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.98%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "20.58"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.01%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.602"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.03%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0232"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.95%  "

# Row 50
$ws.Range("E50").Value = "  +1.10%  "

# Row 51
$ws.Range("E51").Value = "  +5.08%  "
